$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 31   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/4/2024  Through  3/10/2024"

# --- Crime statistics table updates (rows 15-30) ---

# Row 15
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 2
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -50
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("G15").Value = 2
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("H15").Value = -50
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = -50

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 31
$ws.Range("J16").Value = 23
$ws.Range("K16").Value = 34.782608695652
$ws.Range("L16").Value = 19.230769230769
$ws.Range("M16").Value = -26.190476190476
$ws.Range("N16").Value = -74.796747967479

# Row 17
$ws.Range("C17").Value = 1
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 4
$ws.Range("D17").NumberFormat = "#,##0"
$ws.Range("E17").Value = -75
$ws.Range("E17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -45.454545454545
$ws.Range("I17").Value = 26
$ws.Range("J17").Value = 31
$ws.Range("K17").Value = -16.129032258064
$ws.Range("L17").Value = -33.333333333333
$ws.Range("M17").Value = 62.5
$ws.Range("N17").Value = -52.727272727272

# Row 18
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -60
$ws.Range("I18").Value = 13
$ws.Range("J18").Value = 22
$ws.Range("K18").Value = -40.90909090909
$ws.Range("L18").Value = -43.478260869565
$ws.Range("M18").Value = -71.739130434782
$ws.Range("N18").Value = -93.5960591133

# Row 19
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 87.5
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -6.25
$ws.Range("I19").Value = 94
$ws.Range("J19").Value = 110
$ws.Range("K19").Value = -14.545454545454
$ws.Range("L19").Value = 14.634146341463
$ws.Range("M19").Value = 3.296703296703
$ws.Range("N19").Value = -25.984251968503

# Row 20
$ws.Range("C20").Value = 9
$ws.Range("E20").Value = 125
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 27
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = 58.823529411764
$ws.Range("L20").Value = 17.391304347826
$ws.Range("M20").Value = -27.027027027027
$ws.Range("N20").Value = -95.296167247386

# Row 21
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -13.483146067415
$ws.Range("I21").Value = 194
$ws.Range("J21").Value = 208
$ws.Range("K21").Value = -6.730769230769
$ws.Range("L21").Value = -0.51282051282
$ws.Range("M21").Value = -17.094017094017
$ws.Range("N21").Value = -82.283105022831

# Row 23
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C23").NumberFormat = "General"
$ws.Range("D23").Value = 2
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -66.666666666666
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = 75

# Row 24
$ws.Range("C24").Value = 32
$ws.Range("E24").Value = 60
$ws.Range("F24").Value = 109
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = 15.95744680851
$ws.Range("I24").Value = 273
$ws.Range("J24").Value = 215
$ws.Range("K24").Value = 26.976744186046
$ws.Range("L24").Value = 43.684210526315
$ws.Range("M24").Value = 60.588235294117

# Row 25
$ws.Range("C25").Value = 31
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 244.444444444444
$ws.Range("F25").Value = 95
$ws.Range("G25").Value = 53
$ws.Range("H25").Value = 79.245283018867
$ws.Range("I25").Value = 214
$ws.Range("J25").Value = 137
$ws.Range("K25").Value = 56.204379562043
$ws.Range("L25").Value = 78.333333333333

# Row 26
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -28.571428571428
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 12.5
$ws.Range("I26").Value = 58
$ws.Range("J26").Value = 52
$ws.Range("K26").Value = 11.538461538461
$ws.Range("L26").Value = 31.818181818181
$ws.Range("M26").Value = 9.43396226415

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 3
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = -66.666666666666
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = -57.142857142857
$ws.Range("L27").Value = 50

# Row 28
$ws.Range("D28").Value = 1
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -66.666666666666
$ws.Range("J28").Value = 9
$ws.Range("K28").Value = -66.666666666666

# Row 29
$ws.Range("L29").Value = -83.333333333333

# Row 30
$ws.Range("L30").Value = -66.666666666666
